$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: fill in the rest of the row (F53 already had a style, now gets a value) ---

# D53: new address string, same formatting as other address cells in column D (e.g. D14)
$ws.Range("D14").Copy()
$ws.Range("D53").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D53").Value = "종로 신설동역한양립스"

$ws.Range("F53").Value = 5
$ws.Range("G53").Formula = "=3.305785*F53"
$ws.Range("H53").Value = 11
$ws.Range("I53").Value = 0

# J53:N53 and P53:Q53 share the same "checkbox" formatting as the rest of the table (e.g. row 14)
$ws.Range("J14:N14").Copy()
$ws.Range("J53:N53").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P14:Q14").Copy()
$ws.Range("P53:Q53").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J53").Value = 1
$ws.Range("K53").Value = 1
$ws.Range("L53").Value = 1
$ws.Range("M53").Value = 1
$ws.Range("N53").Value = 1
$ws.Range("O53").Formula = "=(P53+Q53)*100+R53"
$ws.Range("P53").Value = 7
$ws.Range("Q53").Value = 60
$ws.Range("R53").Value = 500

# --- Row 54: blank address cell (formatted like D22) plus a lone value in H ---
$ws.Range("D22").Copy()
$ws.Range("D54").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H54").Value = 123

# --- Row 56 (row 55 intentionally left blank) ---
$ws.Range("H56").Value = 4.45

# --- Update the visible selection to match the latest edits ---
$ws.Activate()
$ws.Range("H57").Select()

$excel.CutCopyMode = $false
